$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B17").Value = " / Barcelona"
$ws.Range("B18").Value = " / Real Sociedad"
$ws.Range("B19").Value = " / Rayo Vallecano"
$ws.Range("B20").Value = " / Real Madrid"
$ws.Range("B21").Value = " / Atlético Madrid"
$ws.Range("B22").Value = " / Barcelona / Real Madrid"
$ws.Range("B23").Value = " / Real Sociedad"
